# Add a new worksheet named "self" after the last existing sheet (icd10),
# and populate it with the labels/description/code cross-check table.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "self"

$data = @(
    ,@("labels_new_self", "description_new_self", "code_new_self", "cross_check_self", "description_old_self", "code_old_self")
    ,@("Ischaemic heart disease", "iron deficiency anaemia", "1331", "missing", "Y", "1074")
    ,@("Ischaemic heart disease", "pernicious anaemia", "1332", "missing", "Y", "1075")
    ,@("Ischaemic heart disease", "aplastic anaemia", "1339", "missing", $null, $null)
    ,@("Ischaemic heart disease", "anaemia", "1447", "missing", $null, $null)
    ,@("Hypertensive diseases", "colitis/not crohns or ulcerative colitis", "1460", "missing", "N", "1065")
    ,@("Hypertensive diseases", "inflammatory bowel disease", "1462", "missing", "Y", "1072")
    ,@("Hypertensive diseases", "crohns disease", "1463", "missing", "Y", "1076")
    ,@("Stroke", "colitis/not crohns or ulcerative colitis", "1460", "missing", "Y", "1081")
    ,@("Stroke", "inflammatory bowel disease", "1462", "missing", "Y", "1086")
    ,@("Stroke", "ulcerative colitis", "1464", "missing", "Y", "1491")
    ,@("Stroke", $null, $null, $null, "Y", "1583")
    ,@("COPD", "dementia/alzheimers/cognitive impairment", "1264", "missing", "Y", "1112")
    ,@("COPD", $null, $null, $null, "Y", "1113")
    ,@("COPD", $null, $null, $null, "Y", "1472")
    ,@("CKD", "epilepsy", "1265", "missing", "Y", "1192")
    ,@("CKD", $null, $null, $null, "Y", "1193")
    ,@("CKD", $null, $null, $null, "Y", "1194")
    ,@("Diabetes", "migraine", "1266", "missing", "Y", "1220")
    ,@("Diabetes", "headaches (not migraine)", "1437", "missing", "Y", "1223")
    ,@("Cirrhosis", "multiple sclerosis", "1262", "missing", "Y", "1158")
    ,@("Cirrhosis", "other demyelinating disease (not multiple sclerosis)", "1398", "missing", "Y", "1506")
    ,@("Cirrhosis", $null, $null, $null, "Y", "1604")
    ,@("Osteoarthritis", "parkinsons disease", "1263", "missing", "Y", "1465")
    ,@("Osteoarthritis", "wolff parkinson white / wpw syndrome", "1485", "missing", $null, $null)
    ,@("Osteoporosis", "sleep apnoea", "1124", "missing", "Y", "1309")
    ,@("Dementia", $null, $null, $null, "Y", "1263")
    ,@("Parkinsonism", $null, $null, $null, "Y", "1262")
    ,@("Multiple sclerosis", "anorexia/bulimia/other eating disorder", "1471", "missing", "Y", "1261")
    ,@("Schizophrenia", "anxiety/panic attacks", "1288", "missing", "Y", "1289")
    ,@("Depression", $null, $null, $null, "Y", "1286")
    ,@("Bipolar", "mania/bipolar disorder/manic depression", "1293", "missing", "Y", "1291")
    ,@("Anemia", "iron deficiency anaemia", "1331", $null, $null, $null)
    ,@("Anemia", "pernicious anaemia", "1332", $null, $null, $null)
    ,@("Anemia", "aplastic anaemia", "1339", $null, $null, $null)
    ,@("Anemia", "anaemia", "1447", $null, $null, $null)
    ,@("Crohn disease", "colitis/not crohns or ulcerative colitis", "1460", $null, $null, $null)
    ,@("Crohn disease", "inflammatory bowel disease", "1462", $null, $null, $null)
    ,@("Crohn disease", "crohns disease", "1463", $null, $null, $null)
    ,@("Ulcerative colitis", "colitis/not crohns or ulcerative colitis", "1460", $null, $null, $null)
    ,@("Ulcerative colitis", "inflammatory bowel disease", "1462", $null, $null, $null)
    ,@("Ulcerative colitis", "ulcerative colitis", "1464", $null, $null, $null)
    ,@("Epilepsy", "epilepsy", "1265", $null, $null, $null)
    ,@("Migraine", "migraine", "1266", $null, $null, $null)
    ,@("Migraine", "headaches (not migraine)", "1437", $null, $null, $null)
    ,@("Sleep apnoea", "sleep apnoea", "1124", $null, $null, $null)
    ,@("Anorexia nervosa", "anorexia/bulimia/other eating disorder", "1471", $null, $null, $null)
    ,@("Anxiety or GAD (not inc. social anxiety)", "anxiety/panic attacks", "1288", $null, $null, $null)
    ,@("Bulimia nervosa", "anorexia/bulimia/other eating disorder", "1471", $null, $null, $null)
    ,@("OCD", "obsessive compulsive disorder (ocd)", "1616", $null, $null, $null)
    ,@("Panic attacks", "anxiety/panic attacks", "1288", $null, $null, $null)
    ,@("Rheumatoid arthritis", "rheumatoid arthritis", "1465", $null, $null, $null)
    ,@("Glomerular diseases", "pyelonephritis", "1516", $null, $null, $null)
    ,@("Glomerular diseases", "nephritis", "1609", $null, $null, $null)
    ,@("Glomerular diseases", "glomerulnephritis", "1610", $null, $null, $null)
    ,@("Renal failure", "renal/kidney failure", "1193", $null, $null, $null)
    ,@("Renal failure", "renal failure requiring dialysis", "1194", $null, $null, $null)
    ,@("Renal failure", "renal failure not requiring dialysis", "1196", $null, $null, $null)
    ,@("Hepatitis", "hepatitis", "1156", $null, $null, $null)
    ,@("Hepatitis", "infective/viral hepatitis", "1157", $null, $null, $null)
    ,@("Hepatitis", "non-infective hepatitis", "1158", $null, $null, $null)
    ,@("Hepatitis", "hepatitis a", "1579", $null, $null, $null)
    ,@("Hepatitis", "hepatitis b", "1580", $null, $null, $null)
    ,@("Hepatitis", "hepatitis c", "1581", $null, $null, $null)
    ,@("Hepatitis", "hepatitis d", "1582", $null, $null, $null)
    ,@("Hepatitis", "hepatitis e", "1583", $null, $null, $null)
    ,@("Asthma", "asthma", "1112", $null, $null, $null)
    ,@("Emphysema", "emphysema/chronic bronchitis", "1114", $null, $null, $null)
    ,@("Emphysema", "emphysema", "1473", $null, $null, $null)
    ,@("Cystic fibrosis", "fibrocystic disease", "1367", $null, $null, $null)
    ,@("Gastro-oesophageal reflux disease (GORD)", "gastro-oesophageal reflux (gord) / gastric reflux", "1139", $null, $null, $null)
    ,@("Oesophagitis", "oesophagitis/barretts oesophagus", "1140", $null, $null, $null)
    ,@("Dermatitis and eczema", "eczema/dermatitis", "1453", $null, $null, $null)
    ,@("Dermatitis and eczema", "contact dermatitis", "1670", $null, $null, $null)
    ,@("Gastritis", "gastritis/gastric erosions", "1154", $null, $null, $null)
)
$numRows = $data.Count
$numCols = 6

# Columns C (code_new_self) and F (code_old_self) hold values that look like
# numbers (e.g. "1331") but must be stored as TEXT, matching the source data.
# Pre-formatting the whole block as text ("@") before writing values forces
# Excel to keep them as strings instead of silently coercing to numbers.
$fullRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($numRows, $numCols))
$fullRange.NumberFormat = "@"

for ($r = 0; $r -lt $numRows; $r++) {
    $rowNum = $r + 1
    $rowData = $data[$r]
    for ($c = 0; $c -lt $numCols; $c++) {
        $val = $rowData[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($rowNum, $c + 1).Value = $val
        }
    }
}

# Restore the default "General" number format / style now that the text
# values are safely stored as shared strings, so the sheet doesn't end up
# with a stray custom number format applied to every cell.
$fullRange.NumberFormat = "General"
$fullRange.Style = "Normal"

# Column widths, matching the authored sheet layout.
$ws.Columns.Item(1).ColumnWidth = 38.7109375
$ws.Columns.Item(2).ColumnWidth = 47.28515625
$ws.Columns.Item(3).ColumnWidth = 14.5703125
$ws.Columns.Item(4).ColumnWidth = 15.85546875
$ws.Columns.Item(5).ColumnWidth = 19.28515625
$ws.Columns.Item(6).ColumnWidth = 13.5703125
